$wb = $excel.ActiveWorkbook

# "Ready for handoff" -> "In Translation" everywhere it appears (Overview!E2:F4
# and the per-locale Status columns zh-cn!C2:C4 / de-de!C2:C4 all share this
# same status text).
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $nRows = $used.Rows.Count
    $nCols = $used.Columns.Count
    for ($r = 1; $r -le $nRows; $r++) {
        for ($c = 1; $c -le $nCols; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            # Cast to [string] (and keep the literal on the left) so PowerShell
            # never coerces the comparison through Boolean for "True"/"False"
            # cell values, which would otherwise match any non-empty string.
            if ("Ready for handoff" -eq [string]$cell.Value2) {
                $cell.Value = "In Translation"
            }
        }
    }
}

# Narrow the "Status" column(s) on every sheet (Overview's zh-cn/de-de status
# columns E & F, and each locale sheet's own Status column C).
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5
